$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Not-Talk"
$ws.Range("F15").Value = "Low Agreeableness"

$ws.Range("A3").Value = 0.6227638125419617
$ws.Range("C3").Value = 2.005167245864868

$ws.Range("A4").Value = 1.501367449760437
$ws.Range("C4").Value = 2.8519487380981445

$ws.Range("A5").Value = 3.0708096027374268
$ws.Range("C5").Value = 5.108794689178467

$ws.Range("A6").Value = 3.308736801147461
$ws.Range("C6").Value = 0.879641056060791

$ws.Range("A7").Value = 4.158206939697266
$ws.Range("C7").Value = 2.8574790954589844

$ws.Range("A8").Value = 4.110445976257324
$ws.Range("C8").Value = 2.7284398078918457

$ws.Range("A9").Value = 2.97580885887146
$ws.Range("C9").Value = 3.5012764930725098
$ws.Range("E9").Value = "None"

$ws.Range("A10").Value = 1.3380838632583618
$ws.Range("C10").Value = 5.163074016571045

$ws.Range("A11").Value = 3.814192771911621
$ws.Range("C11").Value = 8.0220308303833
